# Natmi LR-pairs output following Dr Hou's advice:
# re-ran the cell-cell signalling summary with an added "ECs" cluster, which
# changes the sending/target cluster cross-product from
#   {FAPs, sCs} x {FAPs, M2, sCs}               (6 rows)
# to
#   {ECs, FAPs, sCs} x {ECs, FAPs, M2, sCs}      (12 rows)
# and recomputes every expression/specificity metric for Col11a1-Ddr1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A:T for rows 2:13 -- Sending cluster, Ligand symbol, Receptor
# symbol, Target cluster, then the 16 numeric NATMI metrics.
$data = @(
  @("ECs","Col11a1","Ddr1","ECs",2,0.6666666666666666,0.006469666666666668,0.019409,0.0003189585344045152,0.0003189585344045152,2,0.6666666666666666,0.353194,1.059582,0.01734490358674178,0.01734490358674178,0.002285047448666667,0.020565427038,0.000005532305027414777,0.000005532305027414776),
  @("ECs","Col11a1","Ddr1","FAPs",2,0.6666666666666666,0.006469666666666668,0.019409,0.0003189585344045152,0.0003189585344045152,3,1,1.890286333333333,5.670859,0.09282953335278148,0.09282953335278148,0.01222952248122222,0.110065702331,0.00002960877190765824,0.00002960877190765824),
  @("ECs","Col11a1","Ddr1","M2",2,0.6666666666666666,0.006469666666666668,0.019409,0.0003189585344045152,0.0003189585344045152,3,1,0.2252136666666666,0.6756409999999999,0.01105995383486111,0.01105995383486111,0.001457057352111111,0.013113516169,0.000003527666665748896,0.000003527666665748896),
  @("ECs","Col11a1","Ddr1","sCs",2,0.6666666666666666,0.006469666666666668,0.019409,0.0003189585344045152,0.0003189585344045152,3,1,17.89429033333333,53.68287100000001,0.8787656092256156,0.8787656092256156,0.1157700936932222,1.041930843239,0.0002802897908036933,0.0002802897908036932),
  @("FAPs","Col11a1","Ddr1","ECs",3,1,19.329446,57.988338,0.9529535422244141,0.952953542224414,2,0.6666666666666666,0.353194,1.059582,0.01734490358674178,0.01734490358674178,6.827044350524001,61.443399154716,0.01652888731252652,0.01652888731252652),
  @("FAPs","Col11a1","Ddr1","FAPs",3,1,19.329446,57.988338,0.9529535422244141,0.952953542224414,3,1,1.890286333333333,5.670859,0.09282953335278148,0.09282953335278148,36.53818760470467,328.843688442342,0.0884622326315725,0.08846223263157249),
  @("FAPs","Col11a1","Ddr1","M2",3,1,19.329446,57.988338,0.9529535422244141,0.952953542224414,3,1,0.2252136666666666,0.6756409999999999,0.01105995383486111,0.01105995383486111,4.353255408295333,39.17929867465799,0.01053962218376938,0.01053962218376938),
  @("FAPs","Col11a1","Ddr1","sCs",3,1,19.329446,57.988338,0.9529535422244141,0.952953542224414,3,1,17.89429033333333,53.68287100000001,0.8787656092256156,0.8787656092256156,345.8867187064887,3112.980468358398,0.8374228000965457,0.8374228000965456),
  @("sCs","Col11a1","Ddr1","ECs",3,1,0.9478076666666667,2.843423,0.04672749924118139,0.04672749924118139,2,0.6666666666666666,0.353194,1.059582,0.01734490358674178,0.01734490358674178,0.3347599810206667,3.012839829186,0.0008104839691878408,0.0008104839691878408),
  @("sCs","Col11a1","Ddr1","FAPs",3,1,0.9478076666666667,2.843423,0.04672749924118139,0.04672749924118139,3,1,1.890286333333333,5.670859,0.09282953335278148,0.09282953335278148,1.791627878928556,16.124650910357,0.004337691949301319,0.004337691949301319),
  @("sCs","Col11a1","Ddr1","M2",3,1,0.9478076666666667,2.843423,0.04672749924118139,0.04672749924118139,3,1,0.2252136666666666,0.6756409999999999,0.01105995383486111,0.01105995383486111,0.2134592399047778,1.921133159143,0.0005168039844259736,0.0005168039844259736),
  @("sCs","Col11a1","Ddr1","sCs",3,1,0.9478076666666667,2.843423,0.04672749924118139,0.04672749924118139,3,1,17.89429033333333,53.68287100000001,0.8787656092256156,0.8787656092256156,16.96034556749256,152.643110107433,0.04106251933826626,0.04106251933826626)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}

# Sheet now spans A1:T13 (1 header row + 12 data rows) instead of A1:T7.
Write-Output "Updated Sheet1 rows 2-13 (A1:T13) with the ECs-inclusive NATMI table"